$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-261)
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03).
$ws.Range("C2:C261").Value = 45202
